$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2191535563547262
$ws.Range("C2").Value = 0.5227517458580421
$ws.Range("D2").Value = 0.3961421483377996
$ws.Range("E2").Value = 0.6293982430367911
$ws.Range("F2").Value = 0.6122841308166216
$ws.Range("G2").Value = 14
